$d = $word.ActiveDocument

# Hex color 2C3E50 encoded as Word's BGR integer (B<<16 | G<<8 | R)
$color = 5258796

function Set-Highlight($paraText, $segments) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs($i)
        $t = $p.Range.Text.TrimEnd()
        if ($t -eq $paraText) {
            $scope = $p.Range.Duplicate
            foreach ($seg in $segments) {
                $r = $scope.Duplicate
                $r.Find.Execute($seg) | Out-Null
                $r.Font.Bold = $true
                $r.Font.Color = $color
            }
            return $true
        }
    }
    Write-Output ("Set-Highlight: paragraph not found -> " + $paraText)
    return $false
}

Set-Highlight '• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%' @('23%', '64%') | Out-Null

Set-Highlight '• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%' @('87%', '71%', '±4.2%', '±2.1%') | Out-Null

Set-Highlight '• Wrote RFP and analyzed bids from 1,200 vendors for research platform development' @('1,200') | Out-Null

Set-Highlight '• Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+' @('$400M', '$1B') | Out-Null

Set-Highlight '• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M' @('73.5%', '$4.7M') | Out-Null

Set-Highlight '• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%' @('87%', '71%') | Out-Null
